$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows appended by the R script (rows 130 & 131) ---

# Give the new date cells (column A) the same date/time format as the
# rest of the column (reuse the existing style instead of creating a
# new one, by copying it from the cell directly above).
$ws.Range("A129").Copy() | Out-Null
$ws.Range("A130:A131").PasteSpecial(-4122) | Out-Null

# adj_close (column G) is stored as text in this sheet (it mirrors the
# numeric close price but as a shared string), so force text entry by
# temporarily marking the cells as Text before writing the values.
$ws.Range("G130:G131").NumberFormat = "@"

# Row 130
$ws.Cells.Item(130, 1).Value = 45470.2916666667
$ws.Cells.Item(130, 2).Value = 0
$ws.Cells.Item(130, 3).Value = 3.46000003814697
$ws.Cells.Item(130, 4).Value = 3.46000003814697
$ws.Cells.Item(130, 5).Value = 3.46000003814697
$ws.Cells.Item(130, 6).Value = 3.46000003814697
$ws.Cells.Item(130, 7).Value = "3.46000003814697"
$ws.Cells.Item(130, 8).Value = "AGAIN.MI"

# Row 131
$ws.Cells.Item(131, 1).Value = 45471.6139236111
$ws.Cells.Item(131, 2).Value = 24000
$ws.Cells.Item(131, 3).Value = 3.57999992370605
$ws.Cells.Item(131, 4).Value = 3.29999995231628
$ws.Cells.Item(131, 5).Value = 3.01999998092651
$ws.Cells.Item(131, 6).Value = 3.40000009536743
$ws.Cells.Item(131, 7).Value = "3.40000009536743"
$ws.Cells.Item(131, 8).Value = "AGAIN.MI"

# Restore column G to the default (unstyled) look now that the text
# values are committed, matching the rest of the adj_close column.
$ws.Range("G130:G131").Style = "Normal"
